$p = $ppt.ActivePresentation

# Add a new slide after the existing one, using the "Title and Content"
# layout (ppLayoutText = 2), which matches slideLayout2.xml ("Título e
# conteúdo": a title placeholder + an idx=1 content placeholder).
$s = $p.Slides.Add(2, 2)

# Title placeholder ("Título 1") - two runs, as typed with autocorrect
# turning "2 fase teste " into one run and "git" into a second run.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "2 fase teste "
$title.LanguageID = "pt-BR"
[void]$title.InsertAfter("git")
